# Regenerate merged AHB files
# - Rename the diff-table header captions from the generic "_old"/"_new"
#   suffixes to the concrete format-version tags "_FV2404"/"_FV2410".
# - Turn the data range A1:U93 into a proper Excel Table (ListObject).
# - Freeze the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row captions (row 1, columns A:U).
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Convert the used range into an Excel Table with the renamed headers.
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U93"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3) Freeze the header row (row 1) so it stays pinned while scrolling.
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.SplitRow = 0
